# Structure changes, improve documentation.
#
# 1) HEADER sheet: reorder the A3:B6 label/value pairs -- a new DOMAIN/CATEGORY
#    pair is introduced at the top, and SOURCE_ORG/SOURCE_PERSON move down
#    (their values stay blank).
# 2) REVENUE_ (hidden lookup sheet): a "-" placeholder is inserted at the top
#    of both the A (REVENUE reason codes) and B (REVENUE sub-codes) lookup
#    columns, shifting everything else down by one row.
# 3) EXPENSE_ (hidden lookup sheet): same "-" placeholder shift for column B,
#    while column A gets a brand-new E1..E10 series (replacing the old
#    R1..R5 series which only had 5 entries).
# 4) The data-validation list ranges on REVENUE/EXPENSE that point at the
#    REVENUE_/EXPENSE_ lookup columns are widened to match the new extents.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) HEADER sheet
# ---------------------------------------------------------------------------
$header = $wb.Worksheets.Item("HEADER")

$header.Cells.Item(3, 1).Value = "DOMAIN"
$header.Cells.Item(3, 2).Value = "IFDAT"

$header.Cells.Item(4, 1).Value = "CATEGORY"
$header.Cells.Item(4, 2).Value = "INCOME"

$header.Cells.Item(5, 1).Value = "SOURCE_ORG"
$header.Cells.Item(5, 2).ClearContents()

$header.Cells.Item(6, 1).Value = "SOURCE_PERSON"
$header.Cells.Item(6, 2).ClearContents()

# ---------------------------------------------------------------------------
# 2) REVENUE_ lookup sheet
# ---------------------------------------------------------------------------
$revenueLookup = $wb.Worksheets.Item("REVENUE_")

$revenueA = @("-", "R1", "R2", "R3", "R4", "R5")
for ($i = 0; $i -lt $revenueA.Length; $i++) {
    $revenueLookup.Cells.Item($i + 1, 1).Value = $revenueA[$i]
}

$revenueB = @("-", "A", "B", "D", "H", "M", "N", "Q", "S", "W", "A2", "A3", "A4", "A5", `
    "A10", "A20", "A30", "A_3", "M2", "M_2", "M_3", "W2", "W3", "W4", "W_2", "W_3", `
    "D_2", "H2", "H3", "I", "OA", "OM", "_O", "_U", "_Z")
for ($i = 0; $i -lt $revenueB.Length; $i++) {
    $revenueLookup.Cells.Item($i + 1, 2).Value = $revenueB[$i]
}

# ---------------------------------------------------------------------------
# 3) EXPENSE_ lookup sheet
# ---------------------------------------------------------------------------
$expenseLookup = $wb.Worksheets.Item("EXPENSE_")

$expenseA = @("-", "E1", "E2", "E3", "E4", "E5", "E6", "E7", "E8", "E9", "E10")
for ($i = 0; $i -lt $expenseA.Length; $i++) {
    $expenseLookup.Cells.Item($i + 1, 1).Value = $expenseA[$i]
}

$expenseB = @("-", "A", "B", "D", "H", "M", "N", "Q", "S", "W", "A2", "A3", "A4", "A5", `
    "A10", "A20", "A30", "A_3", "M2", "M_2", "M_3", "W2", "W3", "W4", "W_2", "W_3", `
    "D_2", "H2", "H3", "I", "OA", "OM", "_O", "_U", "_Z")
for ($i = 0; $i -lt $expenseB.Length; $i++) {
    $expenseLookup.Cells.Item($i + 1, 2).Value = $expenseB[$i]
}

# ---------------------------------------------------------------------------
# 4) Widen the data-validation list ranges to cover the new lookup extents
# ---------------------------------------------------------------------------
$revenue = $wb.Worksheets.Item("REVENUE")
$revenue.Range("B4:B20").Validation.Modify(3, 1, 1, "'REVENUE_'!`$A`$1:`$A`$6")
$revenue.Range("D4:D20").Validation.Modify(3, 1, 1, "'REVENUE_'!`$B`$1:`$B`$35")

$expense = $wb.Worksheets.Item("EXPENSE")
$expense.Range("B4:B20").Validation.Modify(3, 1, 1, "'EXPENSE_'!`$A`$1:`$A`$11")
$expense.Range("D4:D20").Validation.Modify(3, 1, 1, "'EXPENSE_'!`$B`$1:`$B`$35")
